# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Updates the "Periodo Mora" values listed in column E (rows 16-20) of the
# account-statement sheet. The periods are re-entered in reverse order
# (2103, 2102, 2101, 2012, 2011) while the middle row (2101) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2103"
$ws.Range("E17").Value = "2102"
$ws.Range("E18").Value = "2101"
$ws.Range("E19").Value = "2012"
$ws.Range("E20").Value = "2011"
